$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("manually_curated_test_set")
Write-Host $ws.Name
